$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Emoji / color-word refactor: swap the "black/red/orange/green" palette
# for a "blue" palette (closed-book emoji set), matching the new
# synthetic array used by the mail-merge status column.
$ws.Cells.Replace("⬛", "📘", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
$ws.Cells.Replace("🟥", "📕", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
$ws.Cells.Replace("🟧", "📙", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
$ws.Cells.Replace("🟩", "📗", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
$ws.Cells.Replace("noir", "bleu", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
